$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1951219512195122
$ws.Range("C2").Value = 0.5701219512195121
$ws.Range("J2").Value = 0.01524390243902439
$ws.Range("P2").Value = 0.1432926829268293
$ws.Range("S2").Value = 0.07621951219512195
$ws.Range("C3").Value = 0.03553299492385787
$ws.Range("J3").Value = 0.005076142131979695
$ws.Range("P3").Value = 0.7360406091370558
$ws.Range("S3").Value = 0.2233502538071066
$ws.Range("J4").Value = 0.0196078431372549
$ws.Range("P4").Value = 0.6862745098039216
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("B6").Value = 0.05676855895196507
$ws.Range("D6").Value = 0.01310043668122271
$ws.Range("F6").Value = 0.0611353711790393
$ws.Range("J6").Value = 0.2139737991266376
$ws.Range("O6").Value = 0.03493449781659388
$ws.Range("Q6").Value = 0.2052401746724891
$ws.Range("R6").Value = 0.03930131004366812
$ws.Range("S6").Value = 0.3755458515283843
$ws.Range("B7").Value = 0.0861244019138756
$ws.Range("D7").Value = 0.01913875598086124
$ws.Range("F7").Value = 0.06220095693779904
$ws.Range("J7").Value = 0.138755980861244
$ws.Range("O7").Value = 0.03349282296650718
$ws.Range("Q7").Value = 0.1626794258373206
$ws.Range("R7").Value = 0.05741626794258373
$ws.Range("S7").Value = 0.4401913875598086
$ws.Range("B8").Value = 0.1363636363636364
$ws.Range("D8").Value = 0.01674641148325359
$ws.Range("E8").Value = 0.002392344497607655
$ws.Range("F8").Value = 0.04066985645933014
$ws.Range("J8").Value = 0.1028708133971292
$ws.Range("O8").Value = 0.01674641148325359
$ws.Range("Q8").Value = 0.1961722488038277
$ws.Range("R8").Value = 0.0645933014354067
$ws.Range("S8").Value = 0.423444976076555
$ws.Range("B9").Value = 0.1159420289855072
$ws.Range("D9").Value = 0.02415458937198068
$ws.Range("F9").Value = 0.06280193236714976
$ws.Range("J9").Value = 0.1304347826086956
$ws.Range("O9").Value = 0.02898550724637681
$ws.Range("Q9").Value = 0.1400966183574879
$ws.Range("R9").Value = 0.06280193236714976
$ws.Range("S9").Value = 0.4347826086956522
$ws.Range("B10").Value = 0.1282271944922547
$ws.Range("D10").Value = 0.02925989672977625
$ws.Range("E10").Value = 0.0008605851979345956
$ws.Range("F10").Value = 0.08433734939759036
$ws.Range("J10").Value = 0.09982788296041308
$ws.Range("O10").Value = 0.02237521514629948
$ws.Range("Q10").Value = 0.1944922547332186
$ws.Range("R10").Value = 0.04905335628227195
$ws.Range("S10").Value = 0.3915662650602409
$ws.Range("G11").Value = 0.1598746081504702
$ws.Range("J11").Value = 0.08463949843260188
$ws.Range("K11").Value = 0.2068965517241379
$ws.Range("L11").Value = 0.5329153605015674
$ws.Range("S11").Value = 0.01567398119122257
$ws.Range("G12").Value = 0.7542857142857143
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.005714285714285714
$ws.Range("L12").Value = 0.01714285714285714
$ws.Range("S12").Value = 0.02285714285714286
$ws.Range("G13").Value = 0.7804878048780488
$ws.Range("J13").Value = 0.2195121951219512
$ws.Range("F15").Value = 0.02369668246445497
$ws.Range("H15").Value = 0.1232227488151659
$ws.Range("I15").Value = 0.07582938388625593
$ws.Range("J15").Value = 0.3270142180094787
$ws.Range("K15").Value = 0.08530805687203792
$ws.Range("M15").Value = 0.009478672985781991
$ws.Range("O15").Value = 0.04265402843601896
$ws.Range("S15").Value = 0.3127962085308057
$ws.Range("F16").Value = 0.01357466063348416
$ws.Range("H16").Value = 0.1221719457013575
$ws.Range("I16").Value = 0.1131221719457014
$ws.Range("J16").Value = 0.4117647058823529
$ws.Range("K16").Value = 0.1131221719457014
$ws.Range("M16").Value = 0.01357466063348416
$ws.Range("O16").Value = 0.04977375565610859
$ws.Range("S16").Value = 0.16289592760181
$ws.Range("F17").Value = 0.01678657074340528
$ws.Range("H17").Value = 0.1918465227817746
$ws.Range("I17").Value = 0.1079136690647482
$ws.Range("J17").Value = 0.4244604316546763
$ws.Range("K17").Value = 0.08393285371702638
$ws.Range("M17").Value = 0.01918465227817746
$ws.Range("O17").Value = 0.04556354916067146
$ws.Range("S17").Value = 0.1103117505995204
$ws.Range("F18").Value = 0.04201680672268908
$ws.Range("H18").Value = 0.2184873949579832
$ws.Range("I18").Value = 0.1176470588235294
$ws.Range("J18").Value = 0.3529411764705883
$ws.Range("K18").Value = 0.09243697478991597
$ws.Range("M18").Value = 0.02521008403361345
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.09243697478991597
$ws.Range("F19").Value = 0.01974723538704581
$ws.Range("H19").Value = 0.2109004739336493
$ws.Range("I19").Value = 0.08530805687203792
$ws.Range("J19").Value = 0.3625592417061612
$ws.Range("K19").Value = 0.1279620853080569
$ws.Range("M19").Value = 0.02132701421800948
$ws.Range("O19").Value = 0.07109004739336493
$ws.Range("S19").Value = 0.1011058451816746
